$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.376.31'
$ws.Range("E2").Value = '  -2.34%  '
$ws.Range("D3").Value = '3.611.05'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.39'
$ws.Range("E5").Value = '  -2.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.52'
$ws.Range("E6").Value = '  -3.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +2.37%  '
$ws.Range("D8").Value = '3.605.29'
$ws.Range("E8").Value = '  -0.35%  '
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.192'
$ws.Range("E10").Value = '  -5.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.60'
$ws.Range("E11").Value = '  +13.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.611'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '48.14'
$ws.Range("E13").Value = '  -4.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000280'
$ws.Range("E14").Value = '  -2.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '682.07'
$ws.Range("E15").Value = '  -1.98%  '
$ws.Range("D16").Value = '4.193.28'
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.96'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = '3.623.13'
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("D19").Value = '70.376.10'
$ws.Range("E19").Value = '  -2.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.122'
$ws.Range("E20").Value = '  -0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.65'
$ws.Range("E21").Value = '  -4.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.37'
$ws.Range("E22").Value = '  -2.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.930'
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.97'
$ws.Range("E24").Value = '  -5.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.26'
$ws.Range("E25").Value = '  -4.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.89'
$ws.Range("E26").Value = '  -3.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.75'
$ws.Range("E27").Value = '  -4.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.66'
$ws.Range("E29").Value = '  -3.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.21'
$ws.Range("E30").Value = '  -2.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.07'
$ws.Range("E31").Value = '  -1.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.24'
$ws.Range("E32").Value = '  -5.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.45'
$ws.Range("E33").Value = '  +1.32%  '
$ws.Range("E34").Value = '  -5.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.94'
$ws.Range("E35").Value = '  -6.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '571.02'
$ws.Range("E36").Value = '  -2.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.02'
$ws.Range("E37").Value = '  -3.32%  '
$ws.Range("E38").Value = '  -3.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '58.47'
$ws.Range("E39").Value = '  -2.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '3.522.25'
$ws.Range("E41").Value = '  -3.92%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0448'
$ws.Range("E42").Value = '  -3.43%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.344'
$ws.Range("E43").Value = '  -2.03%  '
$ws.Range("E44").Value = '  -3.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '34.03'
$ws.Range("E45").Value = '  -5.76%  '
$ws.Range("D46").Value = '0.0₃0721'
$ws.Range("E46").Value = '  -6.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.96'
$ws.Range("E47").Value = '  +4.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.64'
$ws.Range("E48").Value = '  -5.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.134'
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '135.96'
$ws.Range("E50").Value = '  +2.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.85'
$ws.Range("E51").Value = '  -4.51%  '
